# VandV Comparison.xlsx edit
# - Inserts a new row above row 4 (shifting old rows 4..16 down to 5..17)
# - Writes new "LLNL Only" / "All Others" category labels in column M
# - Updates a few numeric values (row3 J/K, and what becomes row10 K)
# - Merges J1:L1 header cells and re-centers the three header groups

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 - everything below shifts down by one
$ws.Rows.Item(4).Insert()

# --- New row 4 content ---
$ws.Range("J4").Value = 1.43
$ws.Range("K4").Value = 0.53
$ws.Range("M4").Value = "LLNL Only"

# --- Row 3 ("HGL Temperature Rise") RP-all values updated ---
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 0.47
$ws.Range("M3").Value = "All Others"

# --- Row 10 ("Carbon Dioxide Concentration") RP-all sigma M value updated ---
$ws.Range("K10").Value = 0.48

# --- Header row (row 1) formatting: center/merge the "RP all" header group ---
$ws.Range("J1:L1").Merge()
$ws.Range("B1:L1").HorizontalAlignment = -4108
$ws.Range("B1:L1").VerticalAlignment = -4108
